# Time vector extraction hack
# - Assign the newly-introduced custom databook pages (Active TB Testing and
#   Treatment / Active TB Progression Rates / Active TB Death Rates) to the
#   relevant "Parameters" rows so the TB framework no longer depends on a
#   single catch-all "parameters" databook page for time-vector extraction.
# - Update the saved UI view state (active sheet/tab, selections, scroll
#   position) to match where the author ended up after making the edits.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Custom Databook Pages" sheet - just a view-state change: the user's
#    selection moved down to A10 (the "Active TB Progression Rates" row).
# ---------------------------------------------------------------------
$wsPages = $wb.Worksheets.Item("Custom Databook Pages")
$wsPages.Activate()
$wsPages.Range("A10").Select()

# ---------------------------------------------------------------------
# 2. "Transitions" sheet - no longer the tab that's active/selected once
#    we're done; the frozen pane's scroll position also moved from column R
#    over to column I. (Re-select the original AJ22:AK22 range afterwards so
#    the active-cell/selection recorded for the sheet is left as it was.)
# ---------------------------------------------------------------------
$wsTrans = $wb.Worksheets.Item("Transitions")
$wsTrans.Activate()
$winTrans = $excel.ActiveWindow
$winTrans.FreezePanes = $false
$wsTrans.Range("I1").Select()
$winTrans.FreezePanes = $true
$wsTrans.Range("AJ22:AK22").Select()

# ---------------------------------------------------------------------
# 3. "Parameters" sheet - the actual data fix: populate column F (Databook
#    Page) for the treatment/progression/death parameter rows that were
#    previously blank, now that their matching custom databook pages exist.
# ---------------------------------------------------------------------
$wsParams = $wb.Worksheets.Item("Parameters")

$wsParams.Range("F51:F74").Value = "sh_atreat"
$wsParams.Range("F75:F84").Value = "sh_aprog"
$wsParams.Range("F85:F96").Value = "sh_death"

# This becomes the active/visible sheet & tab once the edits are done, with
# the view scrolled further down and the new block selected.
$wsParams.Activate()
$wsParams.Range("F75:F84").Select()

Write-Host "Applied time vector extraction hack edits."
